$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 36 and 37 used to hold the "THANKS" / "SUCCESS" dictionary keys
# (both pointing at the shared "Lore Ipsum" placeholder text). They are
# replaced with the new "STOP_HEAD" / "STOP_TEXT" dictionary entries and
# their German / English translations. Row 38 ("AGAIN") stays as-is.
$ws.Range("A36").Value() = "STOP_HEAD"
$ws.Range("B36").Value() = "Entschuldigung!"
$ws.Range("C36").Value() = "Sorry!"

$ws.Range("A37").Value() = "STOP_TEXT"
$ws.Range("B37").Value() = "Ihr Abhörgerät ist leider nicht geeignet für diesen Test."
$ws.Range("C37").Value() = "Your listening device is not suitable for this test."

# Move the saved selection/scroll position down to the newly edited area.
$win = $excel.ActiveWindow
$win.ScrollRow = 31
$win.ScrollColumn = 1
$ws.Range("B37").Select()
